# Updates the "Price" (D) and "Volume(1h)" (E) columns to the freshly scraped
# figures, and swaps the FLOKI / dogwifhat rows (46 <-> 47).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.341.48"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "3.684.74"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "678.44"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.17"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.08"
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("E12").Value = "  -3.28%  "
$ws.Range("D13").Value = "4.308.27"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.38"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "3.693.29"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "69.284.38"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.03"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.61"
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.89"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "3.832.22"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -5.50%  "
$ws.Range("E27").Value = "  -4.39%  "
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.77"
$ws.Range("E30").Value = "  -3.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("D35").Value = "3.674.42"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("E36").Value = "  -5.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.26"
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.24"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0906"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "170.48"
$ws.Range("E43").Value = "  +4.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.943"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.70"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.12"
$ws.Range("E48").Value = "  -6.11%  "
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("E50").Value = "  -5.72%  "
$ws.Range("E51").Value = "  -2.69%  "

# Row 46 (was FLOKI) is now dogwifhat; row 47 (was dogwifhat) is now FLOKI.
# Each coin keeps its own Price, stored as text for the same reason as above,
# and picks up a new Volume(1h) reading for its new row.
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.72"
$ws.Range("E46").Value = "  -4.28%  "
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000280"
$ws.Range("E47").Value = "  -2.65%  "
